# Refresh cached Universalis market-price snapshots + derived profit
# columns (currentAveragePrice*, LevePrice*, LeveProfit*) across all
# eight crafting-leve sheets, per the scheduled market-data pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 4: Root Rush
$ws.Range("H4").Value = 404
$ws.Range("I4").Value = 263.1111
$ws.Range("K4").Value = 263.1111
$ws.Range("M4").Value = -149.1111
# Row 9: Distill, My Heart
$ws.Range("H9").Value = 1994.5714
$ws.Range("J9").Value = 1993.8334
$ws.Range("L9").Value = 1993.8334
$ws.Range("N9").Value = -2331.8334
# Row 40: Stuck in the Moment
$ws.Range("H40").Value = 2415
$ws.Range("I40").Value = 1872.5
$ws.Range("K40").Value = 1872.5
$ws.Range("M40").Value = -1697.5
# Row 55: A Real Smooth Move
$ws.Range("H55").Value = 307.5
$ws.Range("I55").Value = 229.66667
$ws.Range("K55").Value = 229.66667
$ws.Range("M55").Value = -15.66667000000001
# Row 80: Cleansing the Wicked Humours
$ws.Range("H80").Value = 1459.5714
$ws.Range("I80").Value = 539.8
$ws.Range("K80").Value = 1619.4
$ws.Range("M80").Value = -621.3999999999999
# Row 83: Washing Away the Sins (L)
$ws.Range("H83").Value = 1459.5714
$ws.Range("I83").Value = 539.8
$ws.Range("K83").Value = 4858.2
$ws.Range("M83").Value = 133.8000000000002
# Row 138: All-night Crafting
$ws.Range("H138").Value = 3372.5
$ws.Range("I138").Value = 3156
$ws.Range("J138").Value = 3733.3333
$ws.Range("K138").Value = 9468
$ws.Range("L138").Value = 11199.9999
$ws.Range("M138").Value = -4328
$ws.Range("N138").Value = -21479.9999

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 6880.0454
$ws.Range("I32").Value = 4786
$ws.Range("J32").Value = 13999.8
$ws.Range("K32").Value = 4786
$ws.Range("L32").Value = 13999.8
$ws.Range("M32").Value = -4499
$ws.Range("N32").Value = -14573.8
# Row 45: Hollow Hallmarks
$ws.Range("H45").Value = 5257.143
$ws.Range("I45").Value = 1800
$ws.Range("J45").Value = 5833.3335
$ws.Range("K45").Value = 1800
$ws.Range("L45").Value = 5833.3335
$ws.Range("M45").Value = -1423
$ws.Range("N45").Value = -6587.3335
# Row 63: Rivets Run through It
$ws.Range("H63").Value = 4379.636
$ws.Range("J63").Value = 5283.5713
$ws.Range("L63").Value = 5283.5713
$ws.Range("N63").Value = -6655.5713
# Row 66: A Riveting Revival (L)
$ws.Range("H66").Value = 4379.636
$ws.Range("J66").Value = 5283.5713
$ws.Range("L66").Value = 26417.8565
$ws.Range("N66").Value = -33281.85649999999
# Row 97: Ore for Me
$ws.Range("H97").Value = 1516.091
$ws.Range("I97").Value = 1076.3334
$ws.Range("K97").Value = 1076.3334
$ws.Range("M97").Value = -580.3334
# Row 110: Scheduled Maintenance
$ws.Range("H110").Value = 1059.1111
$ws.Range("I110").Value = 1059.1111
$ws.Range("K110").Value = 1059.1111
$ws.Range("M110").Value = 985.8888999999999

$ws = $wb.Worksheets.Item("BSM")
# Row 86: Through Thick and Thin
$ws.Range("H86").Value = 6230.5835
$ws.Range("J86").Value = 6880
$ws.Range("L86").Value = 6880
$ws.Range("N86").Value = -9126
# Row 89: Piercing Eyes Deserve Piercing Shafts (L)
$ws.Range("H89").Value = 6230.5835
$ws.Range("J89").Value = 6880
$ws.Range("L89").Value = 34400
$ws.Range("N89").Value = -45632
# Row 94: High Steal
$ws.Range("H94").Value = 4136.7334
$ws.Range("I94").Value = 3822.818
$ws.Range("K94").Value = 3822.818
$ws.Range("M94").Value = -3371.818
# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 3869.111
$ws.Range("I134").Value = 3977.75
$ws.Range("K134").Value = 11933.25
$ws.Range("M134").Value = -9398.25

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found
$ws.Range("H31").Value = 2651.2727
$ws.Range("I31").Value = 2762.3333
$ws.Range("J31").Value = 2518
$ws.Range("K31").Value = 2762.3333
$ws.Range("L31").Value = 2518
$ws.Range("M31").Value = -2467.3333
$ws.Range("N31").Value = -3108
# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 2651.2727
$ws.Range("I34").Value = 2762.3333
$ws.Range("J34").Value = 2518
$ws.Range("K34").Value = 2762.3333
$ws.Range("L34").Value = 2518
$ws.Range("M34").Value = -2560.3333
$ws.Range("N34").Value = -2922
# Row 107: Built to Last
$ws.Range("H107").Value = 396.33334
$ws.Range("J107").Value = 447
$ws.Range("L107").Value = 447
$ws.Range("N107").Value = -4287
# Row 122: Timber of Tenkonto
$ws.Range("H122").Value = 1132.6666
$ws.Range("I122").Value = 948.1111
$ws.Range("J122").Value = 1686.3334
$ws.Range("K122").Value = 2844.3333
$ws.Range("L122").Value = 5059.0002
$ws.Range("M122").Value = -394.3332999999998
$ws.Range("N122").Value = -9959.0002
# Row 132: Hull Lotta Damage
$ws.Range("H132").Value = 1738.1333
$ws.Range("I132").Value = 1751.5714
$ws.Range("K132").Value = 5254.7142
$ws.Range("M132").Value = -2724.7142

$ws = $wb.Worksheets.Item("CUL")
# Row 2: Pork Is a Salty Food
$ws.Range("H2").Value = 47.583332
$ws.Range("J2").Value = 92.75
$ws.Range("L2").Value = 556.5
$ws.Range("N2").Value = -782.5
# Row 5: What a Sap
$ws.Range("H5").Value = 845.9
$ws.Range("J5").Value = 933.4286
$ws.Range("L5").Value = 2800.2858
$ws.Range("N5").Value = -3024.2858
# Row 9: Jack of All Plates
$ws.Range("H9").Value = 266
$ws.Range("J9").Value = 100
$ws.Range("L9").Value = 300
$ws.Range("N9").Value = -748
# Row 46: Feeding Frenzy
$ws.Range("H46").Value = 887.5
$ws.Range("I46").Value = 875
$ws.Range("K46").Value = 2625
$ws.Range("M46").Value = -2534
# Row 123: Topping Up the Pot
$ws.Range("H123").Value = 2000
$ws.Range("I123").Value = 2000
$ws.Range("K123").Value = 6000
$ws.Range("M123").Value = -3550
# Row 131: The Mountain Steeped
$ws.Range("H131").Value = 1499.6666
$ws.Range("I131").Value = 1100
$ws.Range("J131").Value = 1699.5
$ws.Range("K131").Value = 3300
$ws.Range("L131").Value = 5098.5
$ws.Range("M131").Value = 1740
$ws.Range("N131").Value = -15178.5
# Row 135: Not-so-secret Ingredient
$ws.Range("H135").Value = 845.9
$ws.Range("J135").Value = 933.4286
$ws.Range("L135").Value = 8400.857399999999
$ws.Range("N135").Value = -13470.8574
# Row 136: Simple Is Hardest
$ws.Range("H136").Value = 4130
$ws.Range("I136").Value = 4130
$ws.Range("K136").Value = 12390
$ws.Range("M136").Value = -7290

$ws = $wb.Worksheets.Item("GSM")
# Row 70: Sky Is the Limit
$ws.Range("H70").Value = 33335666
$ws.Range("J70").Value = 3499
$ws.Range("L70").Value = 3499
$ws.Range("N70").Value = -4039
# Row 73: Hulls of Broken Dreams (L)
$ws.Range("H73").Value = 33335666
$ws.Range("J73").Value = 3499
$ws.Range("L73").Value = 3499
$ws.Range("N73").Value = -5371
# Row 80: Needs More Prayerbell
$ws.Range("H80").Value = 3056.3333
$ws.Range("I80").Value = 2869
$ws.Range("K80").Value = 2869
$ws.Range("M80").Value = -1871
# Row 83: With a Noise That Reaches Heaven (L)
$ws.Range("H83").Value = 3056.3333
$ws.Range("I83").Value = 2869
$ws.Range("K83").Value = 14345
$ws.Range("M83").Value = -9353
# Row 97: If I'd a Koppranickel for Every Time...
$ws.Range("H97").Value = 1032.7
$ws.Range("I97").Value = 653.375
$ws.Range("J97").Value = 2550
$ws.Range("K97").Value = 653.375
$ws.Range("L97").Value = 2550
$ws.Range("M97").Value = -157.375
$ws.Range("N97").Value = -3542
# Row 126: Gold Rush Order
$ws.Range("H126").Value = 142861390
$ws.Range("J126").Value = 5274.75
$ws.Range("L126").Value = 15824.25
$ws.Range("N126").Value = -20764.25
# Row 132: On Board for Lar
$ws.Range("H132").Value = 1645.762
$ws.Range("I132").Value = 1663.421
$ws.Range("J132").Value = 1478
$ws.Range("K132").Value = 4990.263
$ws.Range("L132").Value = 4434
$ws.Range("M132").Value = -2460.263
$ws.Range("N132").Value = -9494
# Row 136: Shiny and Good
$ws.Range("H136").Value = 30326
$ws.Range("J136").Value = 30326
$ws.Range("L136").Value = 90978
$ws.Range("N136").Value = -96078

$ws = $wb.Worksheets.Item("LTW")
# Row 16: Saddle Sore
$ws.Range("H16").Value = 1750
$ws.Range("J16").Value = 1500
$ws.Range("L16").Value = 1500
$ws.Range("N16").Value = -1840
# Row 30: Packing a Punch
$ws.Range("H30").Value = 1180.4
$ws.Range("I30").Value = 1180.4
$ws.Range("K30").Value = 1180.4
$ws.Range("M30").Value = -1072.4
# Row 43: Subordinate Clause
$ws.Range("H43").Value = 29999
$ws.Range("J43").Value = 29999
$ws.Range("L43").Value = 29999
$ws.Range("N43").Value = -30385
# Row 46: Supply Side Logic
$ws.Range("H46").Value = 1711.3334
$ws.Range("I46").Value = 1400
$ws.Range("K46").Value = 1400
$ws.Range("M46").Value = -1212
# Row 61: Spelling Me Softly
$ws.Range("H61").Value = 1654.6666
$ws.Range("I61").Value = 1592
$ws.Range("J61").Value = 1780
$ws.Range("K61").Value = 1592
$ws.Range("L61").Value = 1780
$ws.Range("M61").Value = -1390
$ws.Range("N61").Value = -2184
# Row 113: Peace in Rest
$ws.Range("H113").Value = 1654.6666
$ws.Range("I113").Value = 1592
$ws.Range("J113").Value = 1780
$ws.Range("K113").Value = 1592
$ws.Range("L113").Value = 1780
$ws.Range("M113").Value = 578
$ws.Range("N113").Value = -6120

$ws = $wb.Worksheets.Item("WVR")
# Row 5: Hire in the Blood
$ws.Range("H5").Value = 2000
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()
# Row 81: Where the Dragonflies, the Net Catches
$ws.Range("H81").Value = 5735.909
$ws.Range("I81").Value = 3219.8
$ws.Range("J81").Value = 7832.6665
$ws.Range("K81").Value = 6439.6
$ws.Range("L81").Value = 15665.333
$ws.Range("M81").Value = -5378.6
$ws.Range("N81").Value = -17787.333
# Row 84: To Kill a Dragon on Nameday (L)
$ws.Range("H84").Value = 5735.909
$ws.Range("I84").Value = 3219.8
$ws.Range("J84").Value = 7832.6665
$ws.Range("K84").Value = 32198
$ws.Range("L84").Value = 78326.66500000001
$ws.Range("M84").Value = -26894
$ws.Range("N84").Value = -88934.66500000001
# Row 92: Modest Beginnings
$ws.Range("H92").Value = 47499.5
$ws.Range("J92").Value = 47499.5
$ws.Range("L92").Value = 47499.5
$ws.Range("N92").Value = -52491.5
# Row 107: Flax Wax
$ws.Range("H107").Value = 280.15384
$ws.Range("I107").Value = 307.875
$ws.Range("K107").Value = 923.625
$ws.Range("M107").Value = 996.375
# Row 126: A Polished Purchase
$ws.Range("H126").Value = 2873.2727
$ws.Range("I126").Value = 3025.6
$ws.Range("K126").Value = 9076.799999999999
$ws.Range("M126").Value = -6606.799999999999
